$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.358.48"
$c.ClearFormats()
$ws.Range("E2").Value = "  -2.03%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.794.62"
$c.ClearFormats()
$ws.Range("E3").Value = "  -1.94%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("E5").Value = "  -0.19%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "307.06"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.31%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4500"
$c.ClearFormats()
$ws.Range("E7").Value = "  -1.49%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3596"
$c.ClearFormats()
$ws.Range("E8").Value = "  -2.54%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "45.88"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.07%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07069"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.53%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.8837"
$c.ClearFormats()
$ws.Range("E11").Value = "  +0.78%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07775"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.88%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "19.36"
$c.ClearFormats()
$ws.Range("E13").Value = "  -1.34%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.770.34"
$c.ClearFormats()
$ws.Range("E14").Value = "  -3.19%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.283"
$c.ClearFormats()
$ws.Range("E15").Value = "  -1.01%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "6.336"
$c.ClearFormats()
$ws.Range("E16").Value = "  -0.85%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "84.86"
$c.ClearFormats()
$ws.Range("E17").Value = "  -2.63%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000008494"
$c.ClearFormats()
$ws.Range("E19").Value = "  -2.61%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.ClearFormats()
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("E21").Value = "  -1.59%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "26.383.97"
$c.ClearFormats()
$ws.Range("E22").Value = "  -2.04%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.984"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.53"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.006.31"
$c.ClearFormats()
$ws.Range("E25").Value = "  -2.04%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.969"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.73%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "151.87"
$c.ClearFormats()
$ws.Range("E27").Value = "  +0.40%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.85"
$c.ClearFormats()
$ws.Range("E28").Value = "  -1.96%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.023"
$c.ClearFormats()
$ws.Range("E29").Value = "  +2.72%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "112.15"
$c.ClearFormats()
$ws.Range("E30").Value = "  -1.55%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.865"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.32%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.08681"
$c.ClearFormats()
$ws.Range("E32").Value = "  -1.34%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.051"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.12%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.736"
$c.ClearFormats()
$ws.Range("E34").Value = "  +6.41%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.440"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.96%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7237"
$c.ClearFormats()
$ws.Range("E36").Value = "  -4.11%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.104"
$c.ClearFormats()
$ws.Range("E37").Value = "  -2.51%  "

$ws.Range("E38").Value = "  +0.14%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.068"
$c.ClearFormats()
$ws.Range("E39").Value = "  -1.99%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01929"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.39%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.05090"
$c.ClearFormats()
$ws.Range("E41").Value = "  -0.99%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.865"
$c.ClearFormats()
$ws.Range("E42").Value = "  -1.04%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.884"
$c.ClearFormats()
$ws.Range("E43").Value = "  -0.64%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.5048"
$c.ClearFormats()
$ws.Range("E44").Value = "  +1.54%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1509"
$c.ClearFormats()
$ws.Range("E45").Value = "  -5.74%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "7.995"
$c.ClearFormats()
$ws.Range("E46").Value = "  -3.74%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.005"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.24%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.4615"
$c.ClearFormats()
$ws.Range("E48").Value = "  -1.46%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "101.10"
$c.ClearFormats()
$ws.Range("E49").Value = "  -1.11%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "9.881"
$c.ClearFormats()
$ws.Range("E50").Value = "  -2.55%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.580"
$c.ClearFormats()
$ws.Range("E51").Value = "  -2.05%  "
